# Fruta / hortaliza, semanal
# Inserts 3 new weekly price rows for "Vega Modelo de Temuco - Melón" just
# above the existing row 477 (pushing the old rows 477-493 down to 480-496),
# then fills the 3 newly-opened rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 477..493 down by 3 to make room for the 3 new records.
$ws.Range("A477:R479").EntireRow.Insert()

function Set-Row {
    param(
        [int]$Row,
        [double]$MercadoId,
        [string]$Mercado,
        [string]$Region,
        [double]$Fecha,
        [double]$Codreg,
        [double]$CategoriaId,
        [string]$Categoria,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMinimo,
        [double]$PrecioMaximo,
        [double]$PrecioPromedio,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgOUnidades,
        [string]$Clasificacion
    )

    $ws.Cells.Item($Row, 1).Value = $MercadoId
    $ws.Cells.Item($Row, 2).Value = $Mercado
    $ws.Cells.Item($Row, 3).Value = $Region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $Codreg
    $ws.Cells.Item($Row, 6).Value = $CategoriaId
    $ws.Cells.Item($Row, 7).Value = $Categoria
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgOUnidades
    $ws.Cells.Item($Row, 18).Value = $Clasificacion
}

# New row 477: Melón / Calameño / Primera
Set-Row 477 10 "Vega Modelo de Temuco" "La Araucanía" `
    44578 9 100112027 "Melón" `
    "Calameño" "Primera" 20000 `
    1000 1000 1000 `
    "`$/unidad" "Región del Maule" 1000 1 `
    "Hortaliza"

# New row 478: Melón / Tuna / Extra
Set-Row 478 10 "Vega Modelo de Temuco" "La Araucanía" `
    44578 9 100112027 "Melón" `
    "Tuna" "Extra" 1000 `
    1300 1300 1300 `
    "`$/unidad" "Región del Maule" 1300 1 `
    "Hortaliza"

# New row 479: Melón / Tuna / Primera
Set-Row 479 10 "Vega Modelo de Temuco" "La Araucanía" `
    44578 9 100112027 "Melón" `
    "Tuna" "Primera" 8000 `
    1000 1000 1000 `
    "`$/unidad" "Región del Maule" 1000 1 `
    "Hortaliza"
